# "Generate Report for Archive"
#
# The localization-status report is regenerated: every cell whose status
# was "Ready for handoff" moves to "In Translation", and the now-narrower
# Status column is re-sized to fit the shorter text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: Status columns are E (zh-cn) and F (de-de) ---
$overview = $wb.Worksheets.Item("Overview")
foreach ($addr in @("E2", "F2", "E3", "F3")) {
    $cell = $overview.Range($addr)
    if ($cell.Text -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}

# --- Per-locale sheets: Status is column C ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in @("C2", "C3")) {
        $cell = $ws.Range($addr)
        if ($cell.Text -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

# Narrow the Status columns to match the shorter "In Translation" text.
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5

$wb.Worksheets.Item("zh-cn").Range("C1").ColumnWidth = 12.5
$wb.Worksheets.Item("de-de").Range("C1").ColumnWidth = 12.5
